$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the B2:H15 data block from numeric values to text strings
# formatted with a comma decimal separator (matching the pasted-in
# source data), except G12 which stays a plain number (7).
$ws.Range("B2").Value = "2,53"
$ws.Range("C2").Value = "2,03"
$ws.Range("D2").Value = "3,28"
$ws.Range("E2").Value = "3,72"
$ws.Range("F2").Value = "2,44"
$ws.Range("G2").Value = "1,69"
$ws.Range("H2").Value = "4,92"
$ws.Range("B3").Value = "3,04"
$ws.Range("C3").Value = "2,48"
$ws.Range("D3").Value = "3,73"
$ws.Range("E3").Value = "4,75"
$ws.Range("F3").Value = "2,83"
$ws.Range("G3").Value = "2,22"
$ws.Range("H3").Value = "5,52"
$ws.Range("B4").Value = "3,49"
$ws.Range("C4").Value = "2,91"
$ws.Range("D4").Value = "4,14"
$ws.Range("E4").Value = "5,19"
$ws.Range("F4").Value = "3,22"
$ws.Range("G4").Value = "2,7"
$ws.Range("H4").Value = "5,99"
$ws.Range("B5").Value = "3,94"
$ws.Range("C5").Value = "3,32"
$ws.Range("D5").Value = "4,58"
$ws.Range("E5").Value = "5,6"
$ws.Range("F5").Value = "3,67"
$ws.Range("G5").Value = "3,15"
$ws.Range("H5").Value = "6,44"
$ws.Range("B6").Value = "4,4"
$ws.Range("C6").Value = "3,78"
$ws.Range("D6").Value = "5,12"
$ws.Range("E6").Value = "6,06"
$ws.Range("F6").Value = "4,08"
$ws.Range("G6").Value = "3,58"
$ws.Range("H6").Value = "6,94"
$ws.Range("B7").Value = "4,85"
$ws.Range("C7").Value = "4,39"
$ws.Range("D7").Value = "5,62"
$ws.Range("E7").Value = "6,55"
$ws.Range("F7").Value = "4,53"
$ws.Range("G7").Value = "4,14"
$ws.Range("H7").Value = "7,42"
$ws.Range("B8").Value = "5,4"
$ws.Range("C8").Value = "4,96"
$ws.Range("D8").Value = "6,13"
$ws.Range("E8").Value = "7,21"
$ws.Range("F8").Value = "5,07"
$ws.Range("G8").Value = "4,62"
$ws.Range("H8").Value = "7,93"
$ws.Range("B9").Value = "5,85"
$ws.Range("C9").Value = "5,52"
$ws.Range("D9").Value = "6,77"
$ws.Range("E9").Value = "7,72"
$ws.Range("F9").Value = "5,74"
$ws.Range("G9").Value = "5,09"
$ws.Range("H9").Value = "8,5"
$ws.Range("B10").Value = "6,54"
$ws.Range("C10").Value = "6,02"
$ws.Range("D10").Value = "7,38"
$ws.Range("E10").Value = "8,28"
$ws.Range("F10").Value = "6,49"
$ws.Range("G10").Value = "5,58"
$ws.Range("H10").Value = "9,14"
$ws.Range("B11").Value = "7,19"
$ws.Range("C11").Value = "6,62"
$ws.Range("D11").Value = "7,9"
$ws.Range("E11").Value = "8,87"
$ws.Range("F11").Value = "7,2"
$ws.Range("G11").Value = "6,24"
$ws.Range("H11").Value = "9,66"
$ws.Range("B12").Value = "7,68"
$ws.Range("C12").Value = "7,14"
$ws.Range("D12").Value = "8,46"
$ws.Range("E12").Value = "9,52"
$ws.Range("F12").Value = "7,87"
$ws.Range("G12").Value = 7
$ws.Range("H12").Value = "10,21"
$ws.Range("B13").Value = "8,21"
$ws.Range("C13").Value = "7,7"
$ws.Range("D13").Value = "9,1"
$ws.Range("E13").Value = "10,16"
$ws.Range("F13").Value = "8,37"
$ws.Range("G13").Value = "7,58"
$ws.Range("H13").Value = "10,8"
$ws.Range("B14").Value = "8,77"
$ws.Range("C14").Value = "8,26"
$ws.Range("D14").Value = "9,63"
$ws.Range("E14").Value = "10,79"
$ws.Range("F14").Value = "8,98"
$ws.Range("G14").Value = "8,2"
$ws.Range("H14").Value = "11,41"
$ws.Range("B15").Value = "9,37"
$ws.Range("C15").Value = "8,83"
$ws.Range("D15").Value = "10,25"
$ws.Range("E15").Value = "11,37"
$ws.Range("F15").Value = "9,45"
$ws.Range("G15").Value = "8,72"
$ws.Range("H15").Value = "11,92"

# Update the active selection to the data block.
$ws.Range("B2:H15").Select()
$excel.ActiveCell = $ws.Range("B2")
